$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add three new header columns (F, G, H), reusing the exact header style
# already present on E1 (bold, centered, thin-bordered) so no new style
# entries are introduced.
$ws.Range("E1").Copy($ws.Range("F1:H1"))

$ws.Range("F1").Value2 = "KNN_Outliers_MAD"
$ws.Range("G1").Value2 = "SVM_Outliers_MAD"
$ws.Range("H1").Value2 = "RF_Outliers_MAD"

# Default all the new "Outliers_MAD" flag cells (rows 2-12) to FALSE.
$ws.Range("F2:H12").Value2 = $false

# The single TRUE outlier flag from the target data: KNN flagged row 6 (Hb 97).
$ws.Range("F6").Value2 = $true

"New dimension: $($ws.UsedRange.Address())"
